$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 currently: 150Ω ±1% 1/8W | R1 | 0805 | C17471
# becomes the merged 470Ω row that also includes R1
$ws.Range("A15").Value = "470Ω ±1% 1/8W"
$ws.Range("B15").Value = "R1 R2 R3 R6 R7"
$ws.Range("D15").Value = "C17710"

# Row 16 currently: 470Ω ±1% 1/8W | R2 R3 R4 R6 R7 | 0805 | C17710
# becomes a new row for R4 at 1.2kΩ
$ws.Range("A16").Value = "1.2kΩ ±1% 1/8W"
$ws.Range("B16").Value = "R4"
$ws.Range("D16").Value = "C17379"
